$d = $word.ActiveDocument

function Invoke-PkgXml($range, [string]$bodyXml) {
    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) "O cancelamento de requisição ocorre conforme:" -> split into several
#    runs, adding "pelo SILOMS Aquisição e Contratos" (italic) in the middle.
# ---------------------------------------------------------------------------
$target = $d.Content.Find.Execute("O cancelamento de requisição ocorre conforme:", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$p1 = $d.Paragraphs.Item(1)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("O cancelamento de requisição")) {
        $p1 = $d.Paragraphs.Item($i)
        break
    }
}
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$body1 = @"
<w:p>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">O cancelamento de requisição </w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">pelo </w:t></w:r>
<w:r><w:rPr><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>SILOMS Aquisição e Contratos</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>ocorre conforme:</w:t></w:r>
</w:p>
"@
Invoke-PkgXml $r1 $body1

# ---------------------------------------------------------------------------
# 2) The two placeholder bullet paragraphs ("ss" / "aa") become real
#    instructions, five brand-new bullet paragraphs are appended, and the
#    final "As requisições canceladas..." paragraph is simplified
#    (also dropping its w:jc="both").
# ---------------------------------------------------------------------------
$pSS = $null
$pNotif = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.TrimEnd([char]13, [char]7) -eq "ss") { $pSS = $d.Paragraphs.Item($i) }
    if ($t.StartsWith("As requisições canceladas")) { $pNotif = $d.Paragraphs.Item($i); break }
}
$r2 = $d.Range($pSS.Range.Start, $pNotif.Range.End - 1)
$body2 = @"
<w:p><w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="120"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>No menu &#8220;Requisição&#8221; c</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">licar </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>em &#8220;Gerenciamento de Requisição&#8221;</w:t></w:r></w:p>
<w:p><w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="120"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>Inserir o número da requisição no campo &#8220;Requisição&#8221;</w:t></w:r></w:p>
<w:p><w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="120"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>Apagar os valores dos campos &#8220;Unidade Requisitante&#8221;, &#8220;Unidade Compradora&#8221; e &#8220;Ano&#8221;</w:t></w:r></w:p>
<w:p><w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="120"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>Clicar no botão com o logo de binóculo.</w:t></w:r></w:p>
<w:p><w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="120"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>Marcar a requisição listada</w:t></w:r></w:p>
<w:p><w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="120"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>Pressionar o botão &#8220;Cancelar Requisição&#8221;</w:t></w:r></w:p>
<w:p><w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="120"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>Clicar no botão &#8220;OK&#8221; no aviso emitido e preencher a justificativa para o cancelamento.</w:t></w:r></w:p>
<w:p><w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="120"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>As requisições canceladas devem ser notificas ao Parque correspondente.</w:t></w:r></w:p>
"@
Invoke-PkgXml $r2 $body2

# ---------------------------------------------------------------------------
# 3) "ATUALIZAR ÍNDICE" heading gains a <w:lastRenderedPageBreak/>.
# ---------------------------------------------------------------------------
$pIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("ATUALIZAR")) { $pIdx = $d.Paragraphs.Item($i); break }
}
$r3 = $d.Range($pIdx.Range.Start, $pIdx.Range.End - 1)
$body3 = @"
<w:p>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="pt-PT"/></w:rPr><w:lastRenderedPageBreak/><w:t>ATUALIZAR ÍNDICE</w:t></w:r>
</w:p>
"@
Invoke-PkgXml $r3 $body3

# ---------------------------------------------------------------------------
# 4) A brand-new page-break-only paragraph is inserted right before the
#    "Elaborado por:" paragraph, which itself gains a
#    <w:lastRenderedPageBreak/>.
# ---------------------------------------------------------------------------
$pElab = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Elaborado por:")) { $pElab = $d.Paragraphs.Item($i); break }
}
$pBefore = $d.Paragraphs.Item($pElab.Range.ParagraphFormat.Parent.Paragraphs.Item(1).Range.Start)
$prevIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $pElab.Range.Start) { $prevIndex = $i - 1; break }
}
$pPrev = $d.Paragraphs.Item($prevIndex)

$r4 = $d.Range($pPrev.Range.End, $pElab.Range.End)
$body4 = @"
<w:p><w:pPr><w:rPr><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br w:type="page"/></w:r></w:p>
<w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Elaborado por:</w:t></w:r></w:p>
"@
Invoke-PkgXml $r4 $body4
